# Dividend Calculation workbook update
# - September (row 11) Taxable Account dividend corrected: 91.54 -> 94.57
# - October (row 12) Taxable Account dividend entered: 0 -> 17.55
# - "All Time" sheet's 2016 total (kept in sync manually) updated to match
# - Selection / active-sheet state updated to match the saved workbook

$wb = $excel.ActiveWorkbook
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# --- Data edits -------------------------------------------------------

# Yearly sheet: September dividend revised upward
$wsYearly.Range("D11").Value = 94.57

# Yearly sheet: October dividend entered (was 0)
$wsYearly.Range("D12").Value = 17.55

# All Time sheet: 2016 taxable-account total kept in sync with Yearly!D15
$wsAllTime.Range("F7").Value = 590.45

# --- View / selection state --------------------------------------------

# "All Time" tab was active before; set its selection/scroll first ...
$wsAllTime.Activate()
$wsAllTime.Range("G40").Select() | Out-Null

# ... then make "Yearly" the active tab with its own selection, matching
# the saved file (tabSelected moves from "All Time" to "Yearly").
$wsYearly.Activate()
$wsYearly.Range("J14").Select() | Out-Null
